# Updated hex values of beq and bne
# Fills in the previously-missing PC, binary and hex values for the
# beq/bne rows (and the PC placeholders that had been left as "??"),
# then clears the red "missing data" highlight from the cells that
# are now filled in, and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- PC column (I): replace the "??" / range placeholders with concrete values ---
$ws.Range("I15").Value = 512
$ws.Range("I16").Value = 648
$ws.Range("I17").Value = 16
$ws.Range("I18").Value = 24
$ws.Range("I20").Value = 32
$ws.Range("I21").Value = 40

# --- beq row (15..18 layout mirrors into row 18's breakdown columns) ---
$ws.Range("M18").Value = "(1000)0000001111101000"
$ws.Range("P18").Value = "0x12AE03E8"

# --- bne row ---
$ws.Range("M20").Value = "(1248)0000010011100000"
$ws.Range("P20").Value = "0x16AE04E0"

# --- Clear the red "incomplete" highlight now that the values are filled in ---
$doneRanges = @("I15", "I16", "I17", "I18", "I19", "I20", "I21", "M18:O18", "M20:O20")
foreach ($addr in $doneRanges) {
    $rng = $ws.Range($addr)
    $rng.Interior.Pattern = -4142
}

# --- Update the view/selection to match the new focal point ---
$ws.Activate()
$ws.Range("P20").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 9
$win.ScrollRow = 1
